$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header change ---
$ws.Range("A1").Value = "Recibo"

# --- Update existing row 2 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "310001"
$ws.Range("A2").ClearFormats()

$ws.Range("C2").Value = "08/07/2024 11:49:45"
$ws.Range("D2").Value = "08/07/2024 11:50:01"
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 1500

# --- Append new rows 3-9 ---
$rows = @(
    @{ r=3; a="310002"; b="AAA01"; c="08/07/2024 18:11:30"; d="09/07/2024 12:12:09"; e="Moto"; f=1500; g=64839; h=14000 },
    @{ r=4; a="310003"; b="BBB01"; c="08/07/2024 18:11:38"; d="09/07/2024 12:13:47"; e="Moto"; f=1500; g=64929; h=14750 },
    @{ r=5; a="310004"; b="CCC01"; c="08/07/2024 18:11:43"; d="09/07/2024 12:14:26"; e="Moto"; f=1500; g=64963; h=14750 },
    @{ r=6; a="310005"; b="DDD01"; c="08/07/2024 18:11:54"; d="09/07/2024 15:13:41"; e="Moto"; f=1500; g=75707; h=14750 },
    @{ r=7; a="310006"; b="EEE01"; c="08/07/2024 18:12:11"; d="09/07/2024 15:14:11"; e="Moto"; f=1500; g=75720; h=14750 },
    @{ r=8; a="310007"; b="FFF01"; c="08/07/2024 18:17:30"; d="09/07/2024 15:15:15"; e="Moto"; f=1500; g=75465; h=15500 },
    @{ r=9; a="310008"; b="GGG01"; c="08/07/2024 18:18:07"; d="09/07/2024 15:15:26"; e="Moto"; f=1500; g=75439; h=15500 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h
}

Write-Output "edit complete"
